$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 151, shifting existing rows 151..254 down to 152..255
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with its data
$ws.Cells.Item(151, 1).Value = 9
$ws.Cells.Item(151, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(151, 3).Value = "Metropolitana"
$ws.Cells.Item(151, 4).Value = 44673
$ws.Cells.Item(151, 5).Value = 13
$ws.Cells.Item(151, 6).Value = 300000001
$ws.Cells.Item(151, 7).Value = "Rabanito"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 5200
$ws.Cells.Item(151, 11).Value = 3000
$ws.Cells.Item(151, 12).Value = 3000
$ws.Cells.Item(151, 13).Value = 3000
$ws.Cells.Item(151, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(151, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(151, 16).Value = 30
$ws.Cells.Item(151, 17).Value = 100
$ws.Cells.Item(151, 18).Value = "Hortaliza"
